$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 41608.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 41608.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 41608.5
$ws.Range("N3").Value = -41836.5

$ws.Range("H32").Value = 1793
$ws.Range("I32").Value = 1937.25
$ws.Range("J32").Value = 1600.6666
$ws.Range("K32").Value = 1937.25
$ws.Range("L32").Value = 1600.6666
$ws.Range("M32").Value = -1611.25
$ws.Range("N32").Value = -2252.6666

$ws.Range("H86").Value = 2623.889
$ws.Range("I86").Value = 2307.647
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 2307.647
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -1184.647
$ws.Range("N86").Value = -10246

$ws.Range("H89").Value = 2623.889
$ws.Range("I89").Value = 2307.647
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 11538.235
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -5922.235000000001
$ws.Range("N89").Value = -51232

$ws.Range("H102").Value = 41608.5
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 41608.5
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 41608.5
$ws.Range("N102").Value = -48098.5

$ws.Range("H112").Value = 1391.5385
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1391.5385
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 4174.6155
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -6390.6155

$ws.Range("H123").Value = 30875
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 30875
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 30875
$ws.Range("N123").Value = -40675

$ws.Range("H124").Value = 34945
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 34945
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 34945
$ws.Range("N124").Value = -44765

$ws.Range("H138").Value = 6099073
$ws.Range("I138").Value = 1271.5238
$ws.Range("J138").Value = 12501765
$ws.Range("K138").Value = 3814.5714
$ws.Range("L138").Value = 37505295
$ws.Range("M138").Value = 1325.4286
$ws.Range("N138").Value = -37515575

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9312.163
$ws.Range("I32").Value = 10916.919
$ws.Range("J32").Value = 4364.1665
$ws.Range("K32").Value = 10916.919
$ws.Range("L32").Value = 4364.1665
$ws.Range("M32").Value = -10629.919

$ws.Range("H122").Value = 5541.1665
$ws.Range("I122").Value = 6189.24
$ws.Range("J122").Value = 2300.8
$ws.Range("K122").Value = 18567.72
$ws.Range("L122").Value = 6902.400000000001
$ws.Range("M122").Value = -16117.72
$ws.Range("N122").Value = -11802.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 39084
$ws.Range("I102").Value = 25556
$ws.Range("J102").Value = 52612
$ws.Range("K102").Value = 25556
$ws.Range("L102").Value = 52612
$ws.Range("M102").Value = -22311

$ws.Range("H107").Value = 2041.1666
$ws.Range("I107").Value = 3026.6
$ws.Range("J107").Value = 1337.2858
$ws.Range("K107").Value = 3026.6
$ws.Range("L107").Value = 1337.2858
$ws.Range("M107").Value = -1106.6
$ws.Range("N107").Value = -5177.2858

$ws.Range("H134").Value = 4105.5884
$ws.Range("I134").Value = 2320.8333
$ws.Range("J134").Value = 8389
$ws.Range("K134").Value = 6962.499899999999
$ws.Range("L134").Value = 25167
$ws.Range("M134").Value = -4427.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 29387.375
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 29387.375
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 29387.375
$ws.Range("N140").Value = -39747.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 700.75
$ws.Range("I47").Value = 700.75
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 2102.25
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -1671.25
$ws.Range("N47").ClearContents()

$ws.Range("H113").Value = 732.6486
$ws.Range("I113").Value = 449.26315
$ws.Range("J113").Value = 1031.7778
$ws.Range("K113").Value = 1347.78945
$ws.Range("L113").Value = 3095.3334
$ws.Range("M113").Value = 822.21055
$ws.Range("N113").Value = -7435.3334

$ws.Range("H134").Value = 4633
$ws.Range("I134").Value = 2976
$ws.Range("J134").Value = 6290
$ws.Range("K134").Value = 8928
$ws.Range("L134").Value = 18870
$ws.Range("M134").Value = -3858
$ws.Range("N134").Value = -29010

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 26063
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 26063
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 26063
$ws.Range("N51").Value = -27081

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H126").Value = 2906.6843
$ws.Range("I126").Value = 1578.24
$ws.Range("J126").Value = 5461.385
$ws.Range("K126").Value = 4734.72
$ws.Range("L126").Value = 16384.155
$ws.Range("M126").Value = -2264.72
$ws.Range("N126").Value = -21324.155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1853.6774
$ws.Range("I16").Value = 1648.8
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 1648.8
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -1478.8

$ws.Range("H30").Value = 20000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 20000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 20000
$ws.Range("N30").Value = -20216

$ws.Range("H40").Value = 6125.4375
$ws.Range("I40").Value = 9633.333000000001
$ws.Range("J40").Value = 4020.7
$ws.Range("K40").Value = 9633.333000000001
$ws.Range("L40").Value = 4020.7
$ws.Range("M40").Value = -9497.333000000001
$ws.Range("N40").Value = -4292.7

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H122").Value = 8700.267
$ws.Range("I122").Value = 10672
$ws.Range("J122").Value = 6975
$ws.Range("K122").Value = 32016
$ws.Range("L122").Value = 20925
$ws.Range("M122").Value = -29566
$ws.Range("N122").Value = -25825

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12855.889
$ws.Range("I62").Value = 4425
$ws.Range("J62").Value = 19600.6
$ws.Range("K62").Value = 4425
$ws.Range("L62").Value = 19600.6
$ws.Range("M62").Value = -3801
$ws.Range("N62").Value = -20848.6

$ws.Range("H65").Value = 12855.889
$ws.Range("I65").Value = 4425
$ws.Range("J65").Value = 19600.6
$ws.Range("K65").Value = 22125
$ws.Range("L65").Value = 98003
$ws.Range("M65").Value = -19005
$ws.Range("N65").Value = -104243
